$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 366.25
$ws.Range("I18").Value = 366.25
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 366.25
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -82.25
$ws.Range("N18").ClearContents()

$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H52").Value = 3216.6667
$ws.Range("I52").Value = 3216.6667
$ws.Range("K52").Value = 9650.000100000001
$ws.Range("M52").Value = -9490.000100000001

$ws.Range("H100").Value = 3244.5625
$ws.Range("I100").Value = 3378.7693
$ws.Range("K100").Value = 3378.7693
$ws.Range("M100").Value = -2837.7693

$ws.Range("H107").Value = 2169
$ws.Range("I107").Value = 2169
$ws.Range("K107").Value = 2169
$ws.Range("M107").Value = -249

$ws.Range("H127").Value = 948.25
$ws.Range("I127").Value = 948.25
$ws.Range("K127").Value = 2844.75
$ws.Range("M127").Value = 2115.25

$ws.Range("H132").Value = 1750.8572
$ws.Range("I132").Value = 1750.8572
$ws.Range("K132").Value = 5252.571599999999
$ws.Range("M132").Value = -2722.571599999999

$ws.Range("H135").Value = 1128.2941
$ws.Range("I135").Value = 886.375
$ws.Range("J135").Value = 4999
$ws.Range("K135").Value = 7977.375
$ws.Range("L135").Value = 44991
$ws.Range("M135").Value = -5442.375
$ws.Range("N135").Value = -50061

$ws.Range("H138").Value = 10748.75
$ws.Range("J138").Value = 10748.75
$ws.Range("L138").Value = 32246.25
$ws.Range("N138").Value = -42526.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12000.2
$ws.Range("I32").Value = 10586.9
$ws.Range("J32").Value = 20480
$ws.Range("K32").Value = 10586.9
$ws.Range("L32").Value = 20480
$ws.Range("M32").Value = -10299.9
$ws.Range("N32").Value = -21054

$ws.Range("H45").Value = 3195
$ws.Range("I45").Value = 2995
$ws.Range("J45").Value = 3995
$ws.Range("K45").Value = 2995
$ws.Range("L45").Value = 3995
$ws.Range("M45").Value = -2618
$ws.Range("N45").Value = -4749

$ws.Range("H61").Value = 350
$ws.Range("I61").Value = 350
$ws.Range("K61").Value = 350
$ws.Range("M61").Value = -138

$ws.Range("H102").Value = 2477.9412
$ws.Range("I102").Value = 2435
$ws.Range("K102").Value = 2435
$ws.Range("M102").Value = -813

$ws.Range("H132").Value = 3459.2
$ws.Range("I132").Value = 2432.6667
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 7298.000100000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -4768.000100000001
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 350
$ws.Range("I136").Value = 350
$ws.Range("K136").Value = 1050
$ws.Range("M136").Value = 1500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8565.75
$ws.Range("I94").Value = 13636.5
$ws.Range("J94").Value = 3495
$ws.Range("K94").Value = 13636.5
$ws.Range("L94").Value = 3495
$ws.Range("M94").Value = -13185.5
$ws.Range("N94").Value = -4397

$ws.Range("H107").Value = 2412.875
$ws.Range("I107").Value = 2216.6
$ws.Range("K107").Value = 2216.6
$ws.Range("M107").Value = -296.5999999999999

$ws.Range("H134").Value = 2451.4
$ws.Range("I134").Value = 2585.6667
$ws.Range("K134").Value = 7757.000100000001
$ws.Range("M134").Value = -5222.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2726.9412
$ws.Range("I31").Value = 2124.818
$ws.Range("K31").Value = 2124.818
$ws.Range("M31").Value = -1829.818

$ws.Range("H34").Value = 2726.9412
$ws.Range("I34").Value = 2124.818
$ws.Range("K34").Value = 2124.818
$ws.Range("M34").Value = -1922.818

$ws.Range("H41").Value = 16999.666
$ws.Range("J41").Value = 24999.5
$ws.Range("L41").Value = 24999.5
$ws.Range("N41").Value = -25855.5

$ws.Range("H51").Value = 49999
$ws.Range("J51").Value = 49999
$ws.Range("L51").Value = 49999
$ws.Range("N51").Value = -51471

$ws.Range("H59").Value = 28987.334
$ws.Range("J59").Value = 28987.334
$ws.Range("L59").Value = 28987.334
$ws.Range("N59").Value = -31277.334

$ws.Range("H61").Value = 49999
$ws.Range("J61").Value = 49999
$ws.Range("L61").Value = 49999
$ws.Range("N61").Value = -50695

$ws.Range("H68").Value = 74295
$ws.Range("J68").Value = 74295
$ws.Range("L68").Value = 74295
$ws.Range("N68").Value = -75793

$ws.Range("H71").Value = 74295
$ws.Range("J71").Value = 74295
$ws.Range("L71").Value = 222885
$ws.Range("N71").Value = -230373

$ws.Range("H74").Value = 71314
$ws.Range("J74").Value = 71314
$ws.Range("L74").Value = 71314
$ws.Range("N74").Value = -73062

$ws.Range("H77").Value = 71314
$ws.Range("J77").Value = 71314
$ws.Range("L77").Value = 213942
$ws.Range("N77").Value = -222678

$ws.Range("H107").Value = 1519.9166
$ws.Range("I107").Value = 1093.4445
$ws.Range("K107").Value = 1093.4445
$ws.Range("M107").Value = 826.5554999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1885.25
$ws.Range("I122").Value = 596.75
$ws.Range("J122").Value = 2529.5
$ws.Range("K122").Value = 5370.75
$ws.Range("L122").Value = 22765.5
$ws.Range("M122").Value = -2920.75
$ws.Range("N122").Value = -27665.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2497.4
$ws.Range("I80").Value = 2122.25
$ws.Range("K80").Value = 2122.25
$ws.Range("M80").Value = -1124.25

$ws.Range("H83").Value = 2497.4
$ws.Range("I83").Value = 2122.25
$ws.Range("K83").Value = 10611.25
$ws.Range("M83").Value = -5619.25

$ws.Range("H132").Value = 4998.5
$ws.Range("I132").Value = 5000
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -12470
$ws.Range("N132").Value = -20054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9490.388999999999
$ws.Range("I22").Value = 10216.1
$ws.Range("J22").Value = 8583.25
$ws.Range("K22").Value = 10216.1
$ws.Range("L22").Value = 8583.25
$ws.Range("M22").Value = -9921.1
$ws.Range("N22").Value = -9173.25

$ws.Range("H24").Value = 27500
$ws.Range("I24").Value = 16250
$ws.Range("J24").Value = 50000
$ws.Range("K24").Value = 16250
$ws.Range("L24").Value = 50000
$ws.Range("M24").Value = -15907
$ws.Range("N24").Value = -50686

$ws.Range("H27").Value = 9490.388999999999
$ws.Range("I27").Value = 10216.1
$ws.Range("J27").Value = 8583.25
$ws.Range("K27").Value = 10216.1
$ws.Range("L27").Value = 8583.25
$ws.Range("M27").Value = -10109.1
$ws.Range("N27").Value = -8797.25

$ws.Range("H55").Value = 732.3333
$ws.Range("I55").Value = 515.1667
$ws.Range("J55").Value = 1166.6666
$ws.Range("K55").Value = 515.1667
$ws.Range("L55").Value = 1166.6666
$ws.Range("M55").Value = -342.1667
$ws.Range("N55").Value = -1512.6666

$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1798

$ws.Range("H82").Value = 2968.6
$ws.Range("I82").Value = 3147.6667
$ws.Range("J82").Value = 2700
$ws.Range("K82").Value = 3147.6667
$ws.Range("L82").Value = 2700
$ws.Range("M82").Value = -2786.6667
$ws.Range("N82").Value = -3422

$ws.Range("H85").Value = 2968.6
$ws.Range("I85").Value = 3147.6667
$ws.Range("J85").Value = 2700
$ws.Range("K85").Value = 3147.6667
$ws.Range("L85").Value = 2700
$ws.Range("M85").Value = -1899.6667
$ws.Range("N85").Value = -5196

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170

$ws.Range("H130").Value = 100000
$ws.Range("J130").Value = 100000
$ws.Range("L130").Value = 100000
$ws.Range("N130").Value = -110040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6949.1
$ws.Range("I62").Value = 4678.2856
$ws.Range("J62").Value = 12247.667
$ws.Range("K62").Value = 4678.2856
$ws.Range("L62").Value = 12247.667
$ws.Range("M62").Value = -4054.2856
$ws.Range("N62").Value = -13495.667

$ws.Range("H65").Value = 6949.1
$ws.Range("I65").Value = 4678.2856
$ws.Range("J65").Value = 12247.667
$ws.Range("K65").Value = 23391.428
$ws.Range("L65").Value = 61238.335
$ws.Range("M65").Value = -20271.428
$ws.Range("N65").Value = -67478.33499999999

$ws.Range("H81").Value = 3024.7
$ws.Range("I81").Value = 2360.7778
$ws.Range("J81").Value = 9000
$ws.Range("K81").Value = 4721.5556
$ws.Range("L81").Value = 18000
$ws.Range("M81").Value = -3660.5556
$ws.Range("N81").Value = -20122

$ws.Range("H84").Value = 3024.7
$ws.Range("I84").Value = 2360.7778
$ws.Range("J84").Value = 9000
$ws.Range("K84").Value = 23607.778
$ws.Range("L84").Value = 90000
$ws.Range("M84").Value = -18303.778
$ws.Range("N84").Value = -100608

$ws.Range("H96").Value = 1966.6666
$ws.Range("I96").Value = 1966.6666
$ws.Range("K96").Value = 1966.6666
$ws.Range("M96").Value = -593.6666

$ws.Range("I113").Value = 1095.3
$ws.Range("J113").Value = 835.6
$ws.Range("K113").Value = 3285.9
$ws.Range("L113").Value = 2506.8
$ws.Range("M113").Value = -1115.9
$ws.Range("N113").Value = -6846.8
